$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.953.16"
$ws.Range("E2").Value = "  +0.71%  "
$ws.Range("D3").Value = "1.749.17"
$ws.Range("E3").Value = "  -0.10%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9966"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.48%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "232.04"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.08%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9978"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.31%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5176"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.34%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2817"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +8.32%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "39.58"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.86%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06126"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.08%  "
$ws.Range("D11").Value = "1.743.16"
$ws.Range("E11").Value = "  -0.42%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.06998"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.86%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "15.49"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.37%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6421"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +6.34%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.527"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.64%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "77.00"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.85%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9964"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.44%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9958"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.52%  "
$ws.Range("D19").Value = "25.949.72"
$ws.Range("E19").Value = "  +0.62%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.53"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.96%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.000006638"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.35%  "
$ws.Range("D22").Value = "1.958.17"
$ws.Range("E22").Value = "  -0.86%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.132"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.561"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.66%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.156"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.17%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "140.04"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.48%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.517"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.62%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.09"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.02%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.823"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.90%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "103.08"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.07%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08309"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.77%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.637"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.81%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.441"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.32%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04414"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.40%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.628"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.82%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9823"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.90%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6092"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.61%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.680"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.61%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01572"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.69%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.933"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.03%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9969"
$ws.Range("D41").Style = "Normal"
$ws.Range("E42").Value = "  -2.38%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.3868"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.73%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.7282"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.25%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.957"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.09%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.05449"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.58%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.376"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +7.53%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1113"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.55%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "52.65"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.50%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "29.84"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.86%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.504"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.67%  "
